$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 4 (Name) - set the Value cell (B4) which was previously empty
$ws.Range("B4").Value = "CategorieprofessionnelleVs"

# Row 8 (Date) - update the date/time value
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
